$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Clear the C4, C6, C8 cell contents (legacy "cuota/period" values removed)
$ws.Range("C4").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("C8").ClearContents()

# Move the small legend table (TC/TARJETA DE CREDITO, C/CREDITO, TD/TARJETA DEBITO)
# from P13:Q15 up one row and two columns left, to N12:O14
$ws.Range("P13:Q15").Copy()
$ws.Range("N12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("N12").Value = $ws.Range("P13").Value()
$ws.Range("O12").Value = $ws.Range("Q13").Value()
$ws.Range("N13").Value = $ws.Range("P14").Value()
$ws.Range("O13").Value = $ws.Range("Q14").Value()
$ws.Range("N14").Value = $ws.Range("P15").Value()
$ws.Range("O14").Value = $ws.Range("Q15").Value()

$ws.Range("P13:Q15").Clear()

# Update view: scroll/selection/zoom
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("Q14").Select()
$ws.Application.ActiveWindow.Zoom = 100

$wb.Save()
